$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8705909252166748
$ws.Range("B1").Value = 2.562511920928955
$ws.Range("C1").Value = 1.313771605491638
$ws.Range("D1").Value = 1.315277457237244
$ws.Range("E1").Value = 1.422910809516907
